$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column (H) mirrors the header formatting used by the other
# header cells (e.g. G1 "sum") -- copy G1's full format (incl. style) onto
# H1 before writing the header text so the style index is reused instead
# of a new one being fabricated.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Data row value for the new column.
$ws.Range("H2").Value = 1
